$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Tax %" column header, matching the style of the other headers in row 1 (M1 = s="3")
$ws.Range("M1").Value = "Tax %"
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

# Add sample tax percentage value under the new header
$ws.Range("M2").Value = 5

# Fix typo in invoice status legend: "Unpai=0" -> "Unpaid=0"
$ws.Range("O5").Value = "Unpaid=0"

# Move selection to O5 (matches saved sheetView selection in the edited workbook)
$ws.Range("O5").Select()
